$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Update formula names for the new "CreateVariableSQL" node (lowercase variant)
$ws.Cells.Item(3,3).Value = "createVariableSQL"
$ws.Cells.Item(4,3).Value = "createVariableSQL"
$ws.Cells.Item(5,3).Value = "createVariableSQL"

# Rename the R-code executor to the SQL-code executor
$ws.Cells.Item(6,3).Value = "executeSQLcode"

# Remove the obsolete rows (delete bottom-up so row numbers stay valid)
# Row 16 holds "validateRCode"
$ws.Rows.Item(16).Delete()
# Row 8 holds "writeButton_onClick"
$ws.Rows.Item(8).Delete()

# Update the selection to match the post-edit state
$ws.Range("C23").Select()
